$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Matheus"
$ws.Range("B1").Value = "teste"
$ws.Range("C1").Value = "tetse "
